# Apply the cryptos-list refresh described in the commit:
#  - updated Price (col D) / Volume(1h) (col E) figures
#  - BitcoinCash/Polkadot (rows 19-20) and dogwifhat/Maker (rows 42-43) swapped places
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "'67.386.67"
$ws.Range("E2").Value = "'  -1.90%  "
$ws.Range("D3").Value = "'3.243.61"
$ws.Range("E3").Value = "'  -5.34%  "
$ws.Range("E4").Value = "'  +0.21%  "
$ws.Range("D5").Value = "'584.74"
$ws.Range("E5").Value = "'  -4.74%  "
$ws.Range("D6").Value = "'145.89"
$ws.Range("E6").Value = "'  -12.51%  "
$ws.Range("E7").Value = "'  +0.14%  "
$ws.Range("D8").Value = "'3.237.45"
$ws.Range("E8").Value = "'  -5.44%  "
$ws.Range("D9").Value = "'0.533"
$ws.Range("E9").Value = "'  -9.97%  "
$ws.Range("D10").Value = "'0.166"
$ws.Range("E10").Value = "'  -14.42%  "
$ws.Range("D11").Value = "'6.68"
$ws.Range("E11").Value = "'  -3.69%  "
$ws.Range("D12").Value = "'0.494"
$ws.Range("E12").Value = "'  -11.86%  "
$ws.Range("D13").Value = "'0.0000241"
$ws.Range("E13").Value = "'  -10.04%  "
$ws.Range("D14").Value = "'37.27"
$ws.Range("E14").Value = "'  -15.09%  "
$ws.Range("D15").Value = "'3.763.17"
$ws.Range("E15").Value = "'  -5.35%  "
$ws.Range("D16").Value = "'67.411.73"
$ws.Range("E16").Value = "'  -2.04%  "
$ws.Range("D17").Value = "'3.250.45"
$ws.Range("E17").Value = "'  -5.16%  "
$ws.Range("E18").Value = "'  -6.39%  "
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").Value = "'6.94"
$ws.Range("E19").Value = "'  -14.27%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "'511.56"
$ws.Range("E20").Value = "'  -11.08%  "
$ws.Range("D21").Value = "'14.61"
$ws.Range("E21").Value = "'  -14.19%  "
$ws.Range("D22").Value = "'0.737"
$ws.Range("E22").Value = "'  -12.26%  "
$ws.Range("D23").Value = "'7.59"
$ws.Range("E23").Value = "'  -15.30%  "
$ws.Range("D24").Value = "'84.16"
$ws.Range("E24").Value = "'  -11.54%  "
$ws.Range("D25").Value = "'13.08"
$ws.Range("E25").Value = "'  -12.55%  "
$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "'  -0.32%  "
$ws.Range("D27").Value = "'3.17"
$ws.Range("E27").Value = "'  -12.40%  "
$ws.Range("D28").Value = "'2.09"
$ws.Range("E28").Value = "'  -12.69%  "
$ws.Range("D29").Value = "'7.80"
$ws.Range("E29").Value = "'  -8.45%  "
$ws.Range("D30").Value = "'28.32"
$ws.Range("E30").Value = "'  -12.83%  "
$ws.Range("D31").Value = "'1.17"
$ws.Range("E31").Value = "'  -4.75%  "
$ws.Range("D32").Value = "'2.58"
$ws.Range("E32").Value = "'  -6.31%  "
$ws.Range("D33").Value = "'6.36"
$ws.Range("E33").Value = "'  -18.31%  "
$ws.Range("E34").Value = "'  +0.12%  "
$ws.Range("D35").Value = "'5.54"
$ws.Range("E35").Value = "'  -15.20%  "
$ws.Range("D36").Value = "'55.57"
$ws.Range("E36").Value = "'  -1.26%  "
$ws.Range("D37").Value = "'503.17"
$ws.Range("E37").Value = "'  -14.41%  "
$ws.Range("D38").Value = "'0.0427"
$ws.Range("E38").Value = "'  -7.57%  "
$ws.Range("D39").Value = "'0.0832"
$ws.Range("E39").Value = "'  -12.45%  "
$ws.Range("D40").Value = "'0.123"
$ws.Range("E40").Value = "'  -11.82%  "
$ws.Range("D41").Value = "'8.73"
$ws.Range("E41").Value = "'  -16.26%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "'2.902.39"
$ws.Range("E42").Value = "'  -9.64%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "'2.67"
$ws.Range("E43").Value = "'  -12.46%  "
$ws.Range("D44").Value = "'0.259"
$ws.Range("E44").Value = "'  -11.33%  "
$ws.Range("E45").Value = "'  -0.07%  "
$ws.Range("D46").Value = "'2.13"
$ws.Range("E46").Value = "'  -9.96%  "
$ws.Range("D47").Value = "'26.02"
$ws.Range("E47").Value = "'  -16.18%  "
$ws.Range("D48").Value = "'0.0₃0551"
$ws.Range("E48").Value = "'  -17.99%  "
$ws.Range("D49").Value = "'123.62"
$ws.Range("E49").Value = "'  -6.69%  "
$ws.Range("E50").Value = "'  -11.21%  "
$ws.Range("D51").Value = "'2.24"
$ws.Range("E51").Value = "'  -18.90%  "
